$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 37

$ws.Cells.Item($row, 1).Value = 10
$ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($row, 3).Value = "La Araucanía"
$ws.Cells.Item($row, 4).Value = 44418
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 5).Value = 9
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100107
$ws.Cells.Item($row, 8).Value = "Otros"
$ws.Cells.Item($row, 9).Value = 100107011
$ws.Cells.Item($row, 10).Value = "Tuna"
$ws.Cells.Item($row, 11).Value = "Sin especificar"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 40
$ws.Cells.Item($row, 14).Value = 30000
$ws.Cells.Item($row, 15).Value = 30000
$ws.Cells.Item($row, 16).Value = 30000
$ws.Cells.Item($row, 17).Value = "$/caja 16 kilos"
$ws.Cells.Item($row, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item($row, 19).Value = 1875
$ws.Cells.Item($row, 20).Value = 16
